$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "55-45=", "29+50=", "23+37=", "77-18=", "23+57=",
    "9+7=", "25+12=", "20+13=", "31-13=", "10+82=",
    "81-18=", "33-15=", "12+50=", "20-16=", "13+61=",
    "28+32=", "91-49=", "43+23=", "0+88=", "25+20=",
    "39-1=", "9+86=", "16+18=", "70-15=", "79-57=",
    "34-18=", "85+2=", "7+42=", "72+26=", "35+3=",
    "56-32=", "50+39=", "33-31=", "85-83=", "75-34=",
    "98-71=", "51-26=", "51-32=", "8+34=", "59+32=",
    "87-58=", "46-26=", "69+7=", "49+13=", "78-28=",
    "40-5=", "12+41=", "60-42=", "65+30=", "59+13=",
    "76+17=", "85-24=", "87-83=", "85-56=", "61-2=",
    "42+29=", "85-10=", "35-25=", "92-56=", "9+45=",
    "97-28=", "93-17=", "58-5=", "76-74=", "92+4=",
    "68-15=", "21+72=", "49+38=", "34-26=", "7+72=",
    "15+70=", "60-21=", "44+5=", "78+12=", "52-28=",
    "94-55=", "41-20=", "51-34=", "72+26=", "66-46=",
    "29-29=", "79-71=", "3+21=", "68-65=", "43-7=",
    "15+38=", "38-33=", "51-49=", "72-39=", "16+44=",
    "72-4=", "71+21=", "79+2=", "76+12=", "45+9=",
    "9+79=", "12+20=", "7+68=", "29+45=", "98-21="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
Write-Host "Done. idx=" $idx
